## case - text area included
## Moves the "Execute" directive from row 5 (TestCaseNumber=...) to row 2
## (now TestCaseNumber=5013), clears row 5's Execute cell (picking up a
## plain/white-fill bordered style instead of the old highlighted one), and
## appends 10 new "Sprint4" test-case rows (298-307) to the Test Cases sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate() | Out-Null

$firstRow = 298
$lastRow = 307

# --- E2: becomes the active "Execute" filter cell ---------------------------
# E5 (untouched so far) already carries the exact target style for E2
# (numFmtId 49 / fillId 4 / borderId 5 / horizontal-left), so copy its
# format across first - this is a single atomic style change (no orphan
# style slots get minted along the way).
$e5Original = $ws.Cells.Item(5, 5)
$e2 = $ws.Cells.Item(2, 5)
$e2.Value = "TestCaseNumber=5013"
$e5Original.Copy() | Out-Null
$e2.PasteSpecial(-4122) | Out-Null

# --- Column A: TestCaseNumber, rows 298-306 (5012-5020) ---------------------
$numbers = @("5012", "5013", "5014", "5015", "5016", "5017", "5018", "5019", "5020")
for ($i = 0; $i -lt $numbers.Length; $i++) {
    $ws.Cells.Item($firstRow + $i, 1).Value = $numbers[$i]
}

# --- Column D: Groups = "Sprint4" (first occurrence, row 298) --------------
$ws.Cells.Item($firstRow, 4).Value = "Sprint4"

# --- Column B: Description, rows 298-307 (Sprint4_Tc_01 .. Sprint4_Tc_10) --
$descriptions = @("Sprint4_Tc_01", "Sprint4_Tc_02", "Sprint4_Tc_03", "Sprint4_Tc_04", "Sprint4_Tc_05", "Sprint4_Tc_06", "Sprint4_Tc_07", "Sprint4_Tc_08", "Sprint4_Tc_09", "Sprint4_Tc_10")
for ($i = 0; $i -lt $descriptions.Length; $i++) {
    $ws.Cells.Item($firstRow + $i, 2).Value = $descriptions[$i]
}

# --- Column D: remaining rows 299-307 reuse "Sprint4" -----------------------
for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = "Sprint4"
}

# --- Column C: Priority, rows 298-307 (1..10, reuse of existing strings) ---
$priorities = @("1", "2", "3", "4", "5", "6", "7", "8", "9", "10")
for ($i = 0; $i -lt $priorities.Length; $i++) {
    $ws.Cells.Item($firstRow + $i, 3).Value = $priorities[$i]
}

# --- Column A: last new TestCaseNumber, row 307 (5021) ----------------------
$ws.Cells.Item($lastRow, 1).Value = "5021"

# Apply the shared "plain text, left aligned" look (style reused from the
# header row's existing text cells) to the whole new block A298:D307 in a
# single formats-only paste per column-block, so no transient/orphan style
# slots get minted.
$plainTextFormatSrc = $ws.Cells.Item($lastRow - 1, 2)   # B306, already style 1
$plainTextFormatSrc.Copy() | Out-Null
$ws.Range("A298:D307").PasteSpecial(-4122) | Out-Null

# --- E5: no longer the active filter - clear it, reset its highlight -------
# E4 already carries the fillId/borderId combo E5 needs (just missing the
# number-format + left alignment), so paste its format in first, then layer
# the two remaining properties directly onto the now-uniquely-owned E5 style
# slot (safe: nothing else references it, so the engine mutates it in place
# instead of minting a throwaway orphan).
$e5 = $ws.Cells.Item(5, 5)
$e4Format = $ws.Cells.Item(4, 5)
$e5.ClearContents() | Out-Null
$e4Format.Copy() | Out-Null
$e5.PasteSpecial(-4122) | Out-Null
$e5.NumberFormat = "@"
$e5.HorizontalAlignment = -4131

# --- Restore the view: scrolled to top, E2 selected -------------------------
$ws.Range("E2").Select() | Out-Null
